$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the "indexVar" worksheet right after "groupDirCols".
# It is a copy of the existing "groupDirRows" template sheet (same layout,
# styles, header row and merged title), with the data-row city cell updated
# to demonstrate the new forEach "indexVar" attribute.
# ---------------------------------------------------------------------------
$groupDirRows = $wb.Worksheets.Item("groupDirRows")
$groupDirCols = $wb.Worksheets.Item("groupDirCols")
$groupDirRows.Copy($null, $groupDirCols)

$indexVarSheet = $wb.Worksheets.Item($groupDirCols.Index + 1)
$indexVarSheet.Name = "indexVar"
$indexVarSheet.Range("A3").Value = "`${index + 1}. `${divisionsList.teams.city}?@indexVar=index"
$indexVarSheet.Activate()
$indexVarSheet.Range("A3").Select()

# ---------------------------------------------------------------------------
# Add the "limit" worksheet right after "indexVar".
# Again a copy of the "groupDirRows" template, demonstrating the new
# forEach "limit" attribute.
# ---------------------------------------------------------------------------
$limitSheet = $wb.Worksheets.Item("groupDirRows")
$limitSheet.Copy($null, $indexVarSheet)

$limitSheet = $wb.Worksheets.Item($indexVarSheet.Index + 1)
$limitSheet.Name = "limit"
$limitSheet.Range("A3").Value = "`${divisionsList.teams.city}?@limit=3"
$limitSheet.Activate()
$limitSheet.Range("A1:E1").Select()
